$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '41.859.85'
$cell.Style = 'Normal'
$ws.Range('E2').Value = '  +0.20%  '

# Row 3
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '2.239.97'
$cell.Style = 'Normal'
$ws.Range('E3').Value = '  +0.76%  '

# Row 4
$ws.Range('E4').Value = '  -1.71%  '

# Row 5
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '250.37'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  +8.76%  '

# Row 6
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '0.634'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  +0.79%  '

# Row 7
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '72.24'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  +6.63%  '

# Row 8
$ws.Range('E8').Value = '  -0.41%  '

# Row 9
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.567'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  +4.28%  '

# Row 10
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '42.21'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  +25.50%  '

# Row 11
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '0.0972'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  +0.95%  '

# Row 12
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '58.52'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  -2.35%  '

# Row 13
$ws.Range('E13').Value = '  +0.62%  '

# Row 14
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '6.92'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  +4.46%  '

# Row 15
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '2.572.27'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -0.06%  '

# Row 16
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '15.13'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  +3.27%  '

# Row 17
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '0.858'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  +1.01%  '

# Row 18
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '2.238.78'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  -0.75%  '

# Row 19
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '41.768.36'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  -0.02%  '

# Row 20
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '0.0₃0968'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  -0.21%  '

# Row 21
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '73.43'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  +0.87%  '

# Row 22
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '6.21'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  +1.28%  '

# Row 23
$ws.Range('B23').Value = 'ImmutableX'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '2.27'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  +29.79%  '

# Row 24
$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '235.63'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  +1.40%  '

# Row 25
$ws.Range('E25').Value = '  +0.72%  '

# Row 26
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '3.73'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  +1.26%  '

# Row 27
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '2.49'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  +6.51%  '

# Row 28
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '10.14'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  +2.77%  '

# Row 29
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '2.20'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  +8.22%  '

# Row 30
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '171.90'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  +4.34%  '

# Row 31
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '20.82'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  +3.74%  '

# Row 32
$ws.Range('E32').Value = '  +1.99%  '

# Row 33
$ws.Range('E33').Value = '  +0.91%  '

# Row 34
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '5.45'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  +5.40%  '

# Row 35
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '0.0723'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  +3.55%  '

# Row 36
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '26.76'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  +35.50%  '

# Row 37
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '4.72'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  +0.81%  '

# Row 38
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '4.12'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  +19.82%  '

# Row 39
$ws.Range('E39').Value = '  +7.17%  '

# Row 40
$ws.Range('E40').Value = '  +4.58%  '

# Row 41
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '6.03'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  +1.26%  '

# Row 42
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '68.40'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  +7.58%  '

# Row 43
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '0.215'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  +18.09%  '

# Row 44
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '5.05'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  +5.89%  '

# Row 45
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '11.60'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  +25.69%  '

# Row 46
$ws.Range('B46').Value = 'SynthetixNetwork'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '4.81'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  +16.64%  '

# Row 47
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '8.72'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  -0.23%  '

# Row 48
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '0.102'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  +1.62%  '

# Row 49
$ws.Range('E49').Value = '  -0.81%  '

# Row 50
$ws.Range('E50').Value = '  +19.31%  '

# Row 51
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '1.19'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  +1.20%  '
